# Trade #17 closed at 2026-02-17 15:18:09 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.83
$wsSummary.Range("B4").Value = -0.17
$wsSummary.Range("B5").Value = -0.2
$wsSummary.Range("B6").Value = 17
$wsSummary.Range("B7").Value = 5
$wsSummary.Range("B9").Value = 29.41

# --- Sheet: Strategy Status (MarketMaking row) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.83
$wsStatus.Range("D4").Value = 17
$wsStatus.Range("E4").Value = -0.17
$wsStatus.Range("F4").Value = -0.17
$wsStatus.Range("G4").Value = 29.41

# --- Sheets: All Trades + MarketMaking (append trade #17 as row 18) ---
$tradeSheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A18").Value = 17

    # Force the date column to remain plain text (matches existing rows) instead
    # of being auto-recognized as a date serial number.
    $ws.Range("B18").NumberFormat = "@"
    $ws.Range("B18").Value = "2026-02-17"
    $ws.Range("B18").Style = "Normal"

    $ws.Range("C18").Value = "15:18:02"
    $ws.Range("D18").Value = "MarketMaking"
    $ws.Range("E18").Value = "DOWN"
    $ws.Range("F18").Value = 0.07000000000000001
    $ws.Range("G18").Value = 0.12
    $ws.Range("H18").Value = "CLOSED"
    $ws.Range("I18").Value = 71.4286
    $ws.Range("J18").Value = 0.05
    $ws.Range("K18").Value = 99.83
    $ws.Range("L18").Value = 0
    $ws.Range("M18").Value = 0
    $ws.Range("N18").Value = 0.6
    $ws.Range("O18").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P18").Value = "early_exit"
    $ws.Range("Q18").Value = 0.15
}
